$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.221.05"
$ws.Range("E2").Value = "  +4.55%  "
$ws.Range("D3").Value = "2.251.73"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.618"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.22%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "2.589.76"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").Value = "2.233.30"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "43.113.71"
$ws.Range("E18").Value = "  +4.73%  "
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +14.64%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +29.77%  "
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0800"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("E36").Value = "  +8.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.66%  "
$ws.Range("E38").Value = "  +18.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.78%  "
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("E42").Value = "  +7.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.85%  "
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.481"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +29.87%  "
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.67%  "
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("D51").Value = "2.462.53"
$ws.Range("E51").Value = "  +3.55%  "
